# Auto-generated edit.ps1 for "Generate Report for Handback"
$wb = $excel.ActiveWorkbook

# ---- Sheet1: Overview ----
$ws1 = $wb.Worksheets.Item("Overview")
while ($ws1.Hyperlinks.Count -gt 0) {
    $__first = $null
    foreach ($__h in $ws1.Hyperlinks) { $__first = $__h; break }
    $__first.Delete()
}
$ws1.Range("A2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.md"
$ws1.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws1.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws1.Range("A3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.md"
$ws1.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws1.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws1.Range("A4").Value2 = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws1.Range("B4").Value2 = "In Translation"
$ws1.Range("C4").Value2 = "In Translation"
$ws1.Range("A5").Value2 = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"
$ws1.Range("A6").Value2 = ".localization-config"
$ws1.Range("B6").Value2 = "Not to be localized"
$ws1.Range("C6").Value2 = "Not to be localized"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/47d403d2-d74c-4b27-a406-98772fcd2329.md", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/777de8af-495f-4a7a-badf-3f57dd19ad29.md", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6a9f09592b78a4c52158fd42a35e01a6847052d1/e2e/c6aa9706-a694-448a-8730-9f92d51da86d.md", [System.Type]::Missing, [System.Type]::Missing, "c6aa9706-a694-448a-8730-9f92d51da86d.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/3f9fb24c-623b-4e25-9556-cc1b30c165f3.md", [System.Type]::Missing, [System.Type]::Missing, "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

# ---- Sheet: zh-cn ----
$ws2 = $wb.Worksheets.Item("zh-cn")
while ($ws2.Hyperlinks.Count -gt 0) {
    $__first = $null
    foreach ($__h in $ws2.Hyperlinks) { $__first = $__h; break }
    $__first.Delete()
}
$ws2.Range("A2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.md"
$ws2.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws2.Range("C2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf"
$ws2.Range("D2").Value2 = "2016-03-10 02:48:12"
$ws2.Range("E2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.md"
$ws2.Range("F2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf"
$ws2.Range("G2").Value2 = "2016-03-10 02:49:02"
$ws2.Range("H2").Value2 = "Include"
$ws2.Range("A3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.md"
$ws2.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws2.Range("C3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf"
$ws2.Range("D3").Value2 = "2016-03-10 02:48:12"
$ws2.Range("E3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.md"
$ws2.Range("F3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf"
$ws2.Range("G3").Value2 = "2016-03-10 02:49:02"
$ws2.Range("H3").Value2 = "Include"
$ws2.Range("A4").Value2 = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws2.Range("B4").Value2 = "In Translation"
$ws2.Range("C4").Value2 = "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.zh-cn.xlf"
$ws2.Range("D4").Value2 = "2016-03-10 02:45:48"
$ws2.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H4").Value2 = "Include"
$ws2.Range("A5").Value2 = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws2.Range("B5").Value2 = "Ready for handoff"
$ws2.Range("C5").Value2 = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.zh-cn.xlf"
$ws2.Range("D5").Value2 = "2016-03-10 02:48:12"
$ws2.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").Value2 = "Include"
$ws2.Range("A6").Value2 = ".localization-config"
$ws2.Range("B6").Value2 = "Not to be localized"
$ws2.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H6").Value2 = "Ignored"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/47d403d2-d74c-4b27-a406-98772fcd2329.md", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd18403d4812d139ed5773c7a14f0d22a952252/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/47d403d2-d74c-4b27-a406-98772fcd2329.md", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd18403d4812d139ed5773c7a14f0d22a952252/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/777de8af-495f-4a7a-badf-3f57dd19ad29.md", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd18403d4812d139ed5773c7a14f0d22a952252/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/777de8af-495f-4a7a-badf-3f57dd19ad29.md", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd18403d4812d139ed5773c7a14f0d22a952252/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6a9f09592b78a4c52158fd42a35e01a6847052d1/e2e/c6aa9706-a694-448a-8730-9f92d51da86d.md", [System.Type]::Missing, [System.Type]::Missing, "c6aa9706-a694-448a-8730-9f92d51da86d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08e0a1a9d82a0448a6836bc130714d3a59baf7f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/3f9fb24c-623b-4e25-9556-cc1b30c165f3.md", [System.Type]::Missing, [System.Type]::Missing, "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfd18403d4812d139ed5773c7a14f0d22a952252/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

# ---- Sheet: de-de ----
$ws3 = $wb.Worksheets.Item("de-de")
while ($ws3.Hyperlinks.Count -gt 0) {
    $__first = $null
    foreach ($__h in $ws3.Hyperlinks) { $__first = $__h; break }
    $__first.Delete()
}
$ws3.Range("A2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.md"
$ws3.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws3.Range("C2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf"
$ws3.Range("D2").Value2 = "2016-03-10 02:48:15"
$ws3.Range("E2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.md"
$ws3.Range("F2").Value2 = "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf"
$ws3.Range("G2").Value2 = "2016-03-10 02:49:07"
$ws3.Range("H2").Value2 = "Include"
$ws3.Range("A3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.md"
$ws3.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws3.Range("C3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf"
$ws3.Range("D3").Value2 = "2016-03-10 02:48:15"
$ws3.Range("E3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.md"
$ws3.Range("F3").Value2 = "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf"
$ws3.Range("G3").Value2 = "2016-03-10 02:49:07"
$ws3.Range("H3").Value2 = "Include"
$ws3.Range("A4").Value2 = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws3.Range("B4").Value2 = "In Translation"
$ws3.Range("C4").Value2 = "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.de-de.xlf"
$ws3.Range("D4").Value2 = "2016-03-10 02:46:10"
$ws3.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H4").Value2 = "Include"
$ws3.Range("A5").Value2 = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws3.Range("B5").Value2 = "Ready for handoff"
$ws3.Range("C5").Value2 = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.de-de.xlf"
$ws3.Range("D5").Value2 = "2016-03-10 02:48:15"
$ws3.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").Value2 = "Include"
$ws3.Range("A6").Value2 = ".localization-config"
$ws3.Range("B6").Value2 = "Not to be localized"
$ws3.Range("D6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("G6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H6").Value2 = "Ignored"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/47d403d2-d74c-4b27-a406-98772fcd2329.md", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c5ba1957ce69c5c1e6ef396a9e2ffcb6ecade8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/47d403d2-d74c-4b27-a406-98772fcd2329.md", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c5ba1957ce69c5c1e6ef396a9e2ffcb6ecade8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "47d403d2-d74c-4b27-a406-98772fcd2329.904cc6c45bd08d3efcd9e0d705af430a165444b7.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/777de8af-495f-4a7a-badf-3f57dd19ad29.md", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c5ba1957ce69c5c1e6ef396a9e2ffcb6ecade8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/777de8af-495f-4a7a-badf-3f57dd19ad29.md", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c5ba1957ce69c5c1e6ef396a9e2ffcb6ecade8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "777de8af-495f-4a7a-badf-3f57dd19ad29.1d0d8c284301540b60aabc3657bcf1e37c8e608e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6a9f09592b78a4c52158fd42a35e01a6847052d1/e2e/c6aa9706-a694-448a-8730-9f92d51da86d.md", [System.Type]::Missing, [System.Type]::Missing, "c6aa9706-a694-448a-8730-9f92d51da86d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1b27aab61a0cdd65fc0cb10b968e58f49867d2a5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/e2e/3f9fb24c-623b-4e25-9556-cc1b30c165f3.md", [System.Type]::Missing, [System.Type]::Missing, "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c5ba1957ce69c5c1e6ef396a9e2ffcb6ecade8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f199f25b84b71b21e8d0904892bf58b131a870de/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

